# Update Leve profit-calculation figures across multiple sheets
# (source data refresh from scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 2252
$ws.Range("I38").Value = 251.875
$ws.Range("J38").Value = 6252.25
$ws.Range("K38").Value = 755.625
$ws.Range("L38").Value = 18756.75
$ws.Range("M38").Value = -383.625
$ws.Range("N38").Value = -19500.75

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1234.1666
$ws.Range("I132").Value = 1301.0667
$ws.Range("J132").Value = 899.6667
$ws.Range("K132").Value = 3903.2001
$ws.Range("L132").Value = 2699.0001
$ws.Range("M132").Value = -1373.2001
$ws.Range("N132").Value = -7759.0001

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6485205.5
$ws.Range("I32").Value = 6701212
$ws.Range("K32").Value = 6701212
$ws.Range("M32").Value = -6700925

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 6716.647
$ws.Range("I61").Value = 5552.75
$ws.Range("K61").Value = 5552.75
$ws.Range("M61").Value = -5340.75

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 13158506
$ws.Range("I110").Value = 13158506
$ws.Range("K110").Value = 13158506
$ws.Range("M110").Value = -13156461

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2086.6
$ws.Range("I122").Value = 2086.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6259.799999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3809.799999999999
$ws.Range("N122").Value = $null

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 6716.647
$ws.Range("I136").Value = 5552.75
$ws.Range("K136").Value = 16658.25
$ws.Range("M136").Value = -14108.25

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 1847.1666
$ws.Range("I22").Value = 1847.1666
$ws.Range("K22").Value = 1847.1666
$ws.Range("M22").Value = -1674.1666

# Row 37: That's Some Fine Grinding
$ws.Range("H37").Value = 1106.8334
$ws.Range("J37").Value = 1517.5
$ws.Range("L37").Value = 1517.5
$ws.Range("N37").Value = -1791.5

# Row 94: High Steal
$ws.Range("H94").Value = 1793.6666
$ws.Range("J94").Value = 1162.5
$ws.Range("L94").Value = 1162.5
$ws.Range("N94").Value = -2064.5

# Row 129: Pruned to Perfection
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

# Row 130: Annals of the Empire I
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 40004540
$ws.Range("I31").Value = 200001070
$ws.Range("K31").Value = 200001070
$ws.Range("M31").Value = -200000775

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 40004540
$ws.Range("I34").Value = 200001070
$ws.Range("K34").Value = 200001070
$ws.Range("M34").Value = -200000868

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 4936.1
$ws.Range("I58").Value = 3091.5
$ws.Range("K58").Value = 3091.5
$ws.Range("M58").Value = -2888.5

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 45995.656
$ws.Range("I132").Value = 3288.4092
$ws.Range("J132").Value = 180218.42
$ws.Range("K132").Value = 9865.2276
$ws.Range("L132").Value = 540655.26
$ws.Range("M132").Value = -7335.2276
$ws.Range("N132").Value = -545715.26

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 7309.9165
$ws.Range("I134").Value = 7651.7896
$ws.Range("J134").Value = 6010.8
$ws.Range("K134").Value = 22955.3688
$ws.Range("L134").Value = 18032.4
$ws.Range("M134").Value = -20420.3688
$ws.Range("N134").Value = -23102.4

# Row 136: Turali Quality
$ws.Range("H136").Value = 4936.1
$ws.Range("I136").Value = 3091.5
$ws.Range("K136").Value = 9274.5
$ws.Range("M136").Value = -6724.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 564.75
$ws.Range("I5").Value = 584.6667
$ws.Range("J5").Value = 505
$ws.Range("K5").Value = 1754.0001
$ws.Range("L5").Value = 1515
$ws.Range("M5").Value = -1642.0001
$ws.Range("N5").Value = -1739

# Row 35: Whirled Peas
$ws.Range("H35").Value = 4673.5
$ws.Range("I35").Value = 647
$ws.Range("J35").Value = 8700
$ws.Range("K35").Value = 1941
$ws.Range("L35").Value = 26100
$ws.Range("M35").Value = -1653
$ws.Range("N35").Value = -26676

# Row 59: Comfort Me with Mushrooms
$ws.Range("H59").Value = 2124.5
$ws.Range("I59").Value = 1749.5
$ws.Range("J59").Value = 2499.5
$ws.Range("K59").Value = 5248.5
$ws.Range("L59").Value = 7498.5
$ws.Range("M59").Value = -4708.5
$ws.Range("N59").Value = -8578.5

# Row 107: Slippery Service
$ws.Range("H107").Value = 1450.963
$ws.Range("I107").Value = 292
$ws.Range("J107").Value = 1714.3636
$ws.Range("K107").Value = 876
$ws.Range("L107").Value = 5143.0908
$ws.Range("M107").Value = 1044
$ws.Range("N107").Value = -8983.0908

# Row 118: Teetotally
$ws.Range("H118").Value = 2500224
$ws.Range("I118").Value = 2500224
$ws.Range("K118").Value = 7500672
$ws.Range("M118").Value = -7499429

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 564.75
$ws.Range("I135").Value = 584.6667
$ws.Range("J135").Value = 505
$ws.Range("K135").Value = 5262.0003
$ws.Range("L135").Value = 4545
$ws.Range("M135").Value = -2727.0003
$ws.Range("N135").Value = -9615

$ws = $wb.Worksheets.Item("GSM")
# Row 36: Keep the Change
$ws.Range("H36").Value = 9466.375
$ws.Range("I36").Value = 20510.334
$ws.Range("J36").Value = 2840
$ws.Range("K36").Value = 20510.334
$ws.Range("L36").Value = 2840
$ws.Range("M36").Value = -20025.334
$ws.Range("N36").Value = -3810

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1271.375
$ws.Range("I97").Value = 1281.5714
$ws.Range("K97").Value = 1281.5714
$ws.Range("M97").Value = -785.5714

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 5189.3335
$ws.Range("I122").Value = 5414.931
$ws.Range("J122").Value = 4535.1
$ws.Range("K122").Value = 16244.793
$ws.Range("L122").Value = 13605.3
$ws.Range("M122").Value = -13794.793
$ws.Range("N122").Value = -18505.3

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 67964.39999999999
$ws.Range("I7").Value = 67964.39999999999
$ws.Range("K7").Value = 67964.39999999999
$ws.Range("M7").Value = -67852.39999999999

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 3148.8333
$ws.Range("I22").Value = 1956.591
$ws.Range("K22").Value = 1956.591
$ws.Range("M22").Value = -1661.591

# Row 27: Fire and Hide
$ws.Range("H27").Value = 3148.8333
$ws.Range("I27").Value = 1956.591
$ws.Range("K27").Value = 1956.591
$ws.Range("M27").Value = -1849.591

# Row 33: Just Rewards
$ws.Range("H33").Value = 1017
$ws.Range("J33").Value = 1017
$ws.Range("L33").Value = 1017
$ws.Range("N33").Value = -1597

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 7800.0713
$ws.Range("I46").Value = 3149.75
$ws.Range("J46").Value = 9660.200000000001
$ws.Range("K46").Value = 3149.75
$ws.Range("L46").Value = 9660.200000000001
$ws.Range("M46").Value = -2961.75
$ws.Range("N46").Value = -10036.2

# Row 126: Battered Books
$ws.Range("H126").Value = 67964.39999999999
$ws.Range("I126").Value = 67964.39999999999
$ws.Range("K126").Value = 203893.2
$ws.Range("M126").Value = -201423.2

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 7284.1
$ws.Range("I132").Value = 6925.15
$ws.Range("K132").Value = 20775.45
$ws.Range("M132").Value = -18245.45

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4462.7144
$ws.Range("I136").Value = 2819.4666
$ws.Range("J136").Value = 6358.769
$ws.Range("K136").Value = 8458.399800000001
$ws.Range("L136").Value = 19076.307
$ws.Range("M136").Value = -5908.399800000001
$ws.Range("N136").Value = -24176.307

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 4528.619
$ws.Range("I126").Value = 4528.619
$ws.Range("K126").Value = 13585.857
$ws.Range("M126").Value = -11115.857

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 7478.9473
$ws.Range("I132").Value = 2334.6667
$ws.Range("J132").Value = 8443.5
$ws.Range("K132").Value = 7004.000100000001
$ws.Range("L132").Value = 25330.5
$ws.Range("M132").Value = -4474.000100000001
$ws.Range("N132").Value = -30390.5
